# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" list (E16:E29) is replaced with a new set of periods
# (2301 down to 2112, i.e. the previous list reversed), and the
# "Salario Basico" figures for the first/last rows (F16 / F29) are swapped
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New "Periodo Mora" values for rows 16-29 (text, column is formatted as @)
$ws.Range("E16").Value = "2301"
$ws.Range("E17").Value = "2212"
$ws.Range("E18").Value = "2211"
$ws.Range("E19").Value = "2210"
$ws.Range("E20").Value = "2209"
$ws.Range("E21").Value = "2208"
$ws.Range("E22").Value = "2207"
$ws.Range("E23").Value = "2206"
$ws.Range("E24").Value = "2205"
$ws.Range("E25").Value = "2204"
$ws.Range("E26").Value = "2203"
$ws.Range("E27").Value = "2202"
$ws.Range("E28").Value = "2201"
$ws.Range("E29").Value = "2112"

# "Salario Basico" values for the first and last row swap places
$ws.Range("F16").Value = 30284
$ws.Range("F29").Value = 36341
